# Updated symbol list on Thu Jan  5 10:00:17 UTC 2023 with GitHub Actions
# Refreshes Price (D), Volume(1h) (E) and Hora (G) columns for every coin row,
# and rotates B/C/D/E for rows 15-19 to reflect the new coin ranking order.
# Values in columns D/E/G are text (e.g. "256.08", "0.47%", "10") rather than
# numbers/percentages, so a leading apostrophe is used to force Excel to store
# them as literal text instead of auto-converting to numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''256.08'
$ws.Range("E2").Value = '''0.47%'
$ws.Range("G2").Value = '''10'
$ws.Range("D3").Value = '''26.84'
$ws.Range("E3").Value = '''-5.02%'
$ws.Range("G3").Value = '''10'
$ws.Range("D4").Value = '''4.717'
$ws.Range("E4").Value = '''-10.35%'
$ws.Range("G4").Value = '''10'
$ws.Range("D5").Value = '''0.05927'
$ws.Range("E5").Value = '''1.29%'
$ws.Range("G5").Value = '''10'
$ws.Range("D6").Value = '''6.660'
$ws.Range("E6").Value = '''-0.78%'
$ws.Range("G6").Value = '''10'
$ws.Range("E7").Value = '''0.18%'
$ws.Range("G7").Value = '''10'
$ws.Range("D8").Value = '''0.9471'
$ws.Range("E8").Value = '''-9.00%'
$ws.Range("G8").Value = '''10'
$ws.Range("D9").Value = '''0.1398'
$ws.Range("E9").Value = '''-0.77%'
$ws.Range("G9").Value = '''10'
$ws.Range("D10").Value = '''0.03844'
$ws.Range("E10").Value = '''11.38%'
$ws.Range("G10").Value = '''10'
$ws.Range("D11").Value = '''0.07111'
$ws.Range("E11").Value = '''-0.91%'
$ws.Range("G11").Value = '''10'
$ws.Range("D12").Value = '''0.03190'
$ws.Range("E12").Value = '''-0.03%'
$ws.Range("G12").Value = '''10'
$ws.Range("D13").Value = '''0.09245'
$ws.Range("E13").Value = '''0.11%'
$ws.Range("G13").Value = '''10'
$ws.Range("D14").Value = '''0.001534'
$ws.Range("E14").Value = '''-1.67%'
$ws.Range("G14").Value = '''10'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '''0.006041'
$ws.Range("E15").Value = '''3.56%'
$ws.Range("G15").Value = '''10'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '''3.483'
$ws.Range("E16").Value = '''-0.43%'
$ws.Range("G16").Value = '''10'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '''3.199'
$ws.Range("E17").Value = '''-0.91%'
$ws.Range("G17").Value = '''10'
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").Value = '''2.219'
$ws.Range("E18").Value = '''-0.16%'
$ws.Range("G18").Value = '''10'
$ws.Range("B19").Value = 'One'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D19").Value = '''0.01041'
$ws.Range("E19").Value = '''1,616.00%'
$ws.Range("G19").Value = '''10'
$ws.Range("E20").Value = '''-1.35%'
$ws.Range("G20").Value = '''10'
$ws.Range("D21").Value = '''0.1299'
$ws.Range("E21").Value = '''-1.32%'
$ws.Range("G21").Value = '''10'
$ws.Range("D22").Value = '''3.812'
$ws.Range("E22").Value = '''7.73%'
$ws.Range("G22").Value = '''10'
$ws.Range("D23").Value = '''0.04215'
$ws.Range("E23").Value = '''1.46%'
$ws.Range("G23").Value = '''10'
$ws.Range("G24").Value = '''10'
$ws.Range("E25").Value = '''-1.31%'
$ws.Range("G25").Value = '''10'
$ws.Range("D26").Value = '''0.004487'
$ws.Range("E26").Value = '''-6.73%'
$ws.Range("G26").Value = '''10'
$ws.Range("D27").Value = '''0.0001200'
$ws.Range("E27").Value = '''-0.08%'
$ws.Range("G27").Value = '''10'
$ws.Range("E28").Value = '''1.89%'
$ws.Range("G28").Value = '''10'
$ws.Range("G29").Value = '''10'
$ws.Range("G30").Value = '''10'
$ws.Range("G31").Value = '''10'
$ws.Range("G32").Value = '''10'
$ws.Range("G33").Value = '''10'
$ws.Range("G34").Value = '''10'
$ws.Range("G35").Value = '''10'
$ws.Range("G36").Value = '''10'
$ws.Range("G37").Value = '''10'
$ws.Range("G38").Value = '''10'
$ws.Range("G39").Value = '''10'
$ws.Range("D40").Value = '''0.03827'
$ws.Range("E40").Value = '''0.62%'
$ws.Range("G40").Value = '''10'
$ws.Range("D41").Value = '''0.006211'
$ws.Range("E41").Value = '''7.35%'
$ws.Range("G41").Value = '''10'
$ws.Range("E42").Value = '''-0.17%'
$ws.Range("G42").Value = '''10'
$ws.Range("D43").Value = '''0.002252'
$ws.Range("E43").Value = '''-4.17%'
$ws.Range("G43").Value = '''10'
$ws.Range("D44").Value = '''0.01056'
$ws.Range("E44").Value = '''8.95%'
$ws.Range("G44").Value = '''10'
$ws.Range("E45").Value = '''5.05%'
$ws.Range("G45").Value = '''10'
$ws.Range("D46").Value = '''0.00000000750'
$ws.Range("E46").Value = '''-0.08%'
$ws.Range("G46").Value = '''10'
$ws.Range("E47").Value = '''-4.89%'
$ws.Range("G47").Value = '''10'
$ws.Range("E48").Value = '''11.17%'
$ws.Range("G48").Value = '''10'
$ws.Range("D49").Value = '''0.00002100'
$ws.Range("E49").Value = '''-0.08%'
$ws.Range("G49").Value = '''10'
$ws.Range("D50").Value = '''0.0002000'
$ws.Range("E50").Value = '''-0.08%'
$ws.Range("G50").Value = '''10'
$ws.Range("G51").Value = '''10'
